# Aula_2/beneficios_consolidado.xlsx
# Remove the "tipo_beneficio_desconto" column (old column F) entirely -
# this shifts percentual_desconto / valor_desconto / valor_liquido one
# column to the left (G->F, H->G, I->H) and also fixes the dimension
# from A1:I6 down to A1:H6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(6).Delete()

# Rows 5 and 6 had an incorrect percentual_desconto (112% / 212%,
# clearly a data-entry bug since it produced a negative valor_liquido).
# Correct it to 5% and recompute the dependent columns.
$ws.Cells.Item(5, 6).Value = 0.05
$ws.Cells.Item(5, 7).Value = 9.26
$ws.Cells.Item(5, 8).Value = 176.03

$ws.Cells.Item(6, 6).Value = 0.05
$ws.Cells.Item(6, 7).Value = 55.66
$ws.Cells.Item(6, 8).Value = 1057.59
